# Hortaliza, Macroferia Regional de Talca - Tomate
# Inserts two new weekly price records (new rows 595 and 596) above the
# existing row 595, shifting the previously-existing rows 595:635 down to
# 597:637 (the used range therefore grows from A1:R635 to A1:R637).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 595 (each Insert() pushes the current
# row 595 and everything below it down by one, so calling it twice yields
# two fresh blank rows at 595/596 and the old row 595 ends up at row 597).
$ws.Rows.Item(595).Insert()
$ws.Rows.Item(595).Insert()

# --- New row 595 ---------------------------------------------------------
$ws.Cells.Item(595, 1).Value = 5
$ws.Cells.Item(595, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(595, 3).Value = "Maule"
$ws.Cells.Item(595, 4).Value = 44714
$ws.Cells.Item(595, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(595, 5).Value = 7
$ws.Cells.Item(595, 6).Value = 100112020
$ws.Cells.Item(595, 7).Value = "Tomate"
$ws.Cells.Item(595, 8).Value = "Larga vida"
$ws.Cells.Item(595, 9).Value = "Primera"
$ws.Cells.Item(595, 10).Value = 1500
$ws.Cells.Item(595, 11).Value = 16000
$ws.Cells.Item(595, 12).Value = 16000
$ws.Cells.Item(595, 13).Value = 16000
$ws.Cells.Item(595, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(595, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(595, 16).Value = 889
$ws.Cells.Item(595, 17).Value = 18
$ws.Cells.Item(595, 18).Value = "Hortaliza"

# --- New row 596 ---------------------------------------------------------
$ws.Cells.Item(596, 1).Value = 5
$ws.Cells.Item(596, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(596, 3).Value = "Maule"
$ws.Cells.Item(596, 4).Value = 44714
$ws.Cells.Item(596, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(596, 5).Value = 7
$ws.Cells.Item(596, 6).Value = 100112020
$ws.Cells.Item(596, 7).Value = "Tomate"
$ws.Cells.Item(596, 8).Value = "Larga vida"
$ws.Cells.Item(596, 9).Value = "Primera"
$ws.Cells.Item(596, 10).Value = 2000
$ws.Cells.Item(596, 11).Value = 8000
$ws.Cells.Item(596, 12).Value = 8000
$ws.Cells.Item(596, 13).Value = 8000
$ws.Cells.Item(596, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(596, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(596, 16).Value = 800
$ws.Cells.Item(596, 17).Value = 10
$ws.Cells.Item(596, 18).Value = "Hortaliza"
